$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 18 previously was the last data row (no border below it). Now that two
# more rows are appended after it, it needs a top+bottom "thin" rule so it
# still reads as the bottom of its own block while being visually separated
# from the newly appended block beneath it.
# Start from the existing bottom-of-block style (row 4's thin-bottom border,
# which already carries the correct wrap/font treatment per column) and then
# add the missing top edge.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A18:E18").Borders.Item(8).LineStyle = 1

# ---------------------------------------------------------------------------
# Row 19: a new dialogue entry re-using the generic "Ah, yes, I've heard!"
# strings already present (same text as rows 2/5/8), just a new line number.
# Copy a full existing row of that shape (values + formats) so the shared
# string references line up automatically, then overwrite the line number.
# ---------------------------------------------------------------------------
$ws.Range("B5:E5").Copy()
$ws.Range("B19:E19").PasteSpecial(-4163)
$ws.Range("B19").Value = 161

# ---------------------------------------------------------------------------
# Row 20: new bottom-of-block entry (same visual treatment as rows 4/7/10/…)
# with brand new text.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("B20").Value = 164
$ws.Range("C20").Value = " You have completed the guild's\ngraduation exam?! How fantastic!"
$ws.Range("D20").Value = " Вы сумели пройти гильдейский\nвыпускной экзамен?! Здорово!"
$ws.Range("E20").Value = " Âú òôíåìé ðñïêóé ãéìûäåêòëéê\nâúðôòëîïê üëèàíåî?! Èäïñïâï!"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet view: the visible window and selection moved down with the new rows.
# ---------------------------------------------------------------------------
$ws.Range("D22").Select()
$ws.Application.ActiveWindow.ScrollRow = 16
